$d = $word.ActiveDocument

$d.Content.Find.Execute("63-26=", $true, $false, $false, $false, $false, $true, 1, $false, "25+62=", 2) | Out-Null
$d.Content.Find.Execute("16+21=", $true, $false, $false, $false, $false, $true, 1, $false, "51-18=", 2) | Out-Null
$d.Content.Find.Execute("69+12=", $true, $false, $false, $false, $false, $true, 1, $false, "58+0=", 2) | Out-Null
$d.Content.Find.Execute("0+89=", $true, $false, $false, $false, $false, $true, 1, $false, "66-65=", 2) | Out-Null
$d.Content.Find.Execute("34-7=", $true, $false, $false, $false, $false, $true, 1, $false, "60+7=", 2) | Out-Null
$d.Content.Find.Execute("51+47=", $true, $false, $false, $false, $false, $true, 1, $false, "57+5=", 2) | Out-Null
$d.Content.Find.Execute("56+23=", $true, $false, $false, $false, $false, $true, 1, $false, "55-1=", 2) | Out-Null
$d.Content.Find.Execute("77-29=", $true, $false, $false, $false, $false, $true, 1, $false, "37+10=", 2) | Out-Null
$d.Content.Find.Execute("5+52=", $true, $false, $false, $false, $false, $true, 1, $false, "25+59=", 2) | Out-Null
$d.Content.Find.Execute("2+28=", $true, $false, $false, $false, $false, $true, 1, $false, "69-61=", 2) | Out-Null
$d.Content.Find.Execute("66-31=", $true, $false, $false, $false, $false, $true, 1, $false, "83+4=", 2) | Out-Null
$d.Content.Find.Execute("31-25=", $true, $false, $false, $false, $false, $true, 1, $false, "3-1=", 2) | Out-Null
$d.Content.Find.Execute("25+61=", $true, $false, $false, $false, $false, $true, 1, $false, "49+20=", 2) | Out-Null
$d.Content.Find.Execute("4+13=", $true, $false, $false, $false, $false, $true, 1, $false, "14+8=", 2) | Out-Null
$d.Content.Find.Execute("72-19=", $true, $false, $false, $false, $false, $true, 1, $false, "36+53=", 2) | Out-Null
$d.Content.Find.Execute("46-42=", $true, $false, $false, $false, $false, $true, 1, $false, "26+45=", 2) | Out-Null
$d.Content.Find.Execute("97-61=", $true, $false, $false, $false, $false, $true, 1, $false, "77-76=", 2) | Out-Null
$d.Content.Find.Execute("14+77=", $true, $false, $false, $false, $false, $true, 1, $false, "66+7=", 2) | Out-Null
$d.Content.Find.Execute("38+9=", $true, $false, $false, $false, $false, $true, 1, $false, "39+9=", 2) | Out-Null
$d.Content.Find.Execute("81-48=", $true, $false, $false, $false, $false, $true, 1, $false, "84-53=", 2) | Out-Null
$d.Content.Find.Execute("55+1=", $true, $false, $false, $false, $false, $true, 1, $false, "64-52=", 2) | Out-Null
$d.Content.Find.Execute("83-16=", $true, $false, $false, $false, $false, $true, 1, $false, "49+48=", 2) | Out-Null
$d.Content.Find.Execute("80-40=", $true, $false, $false, $false, $false, $true, 1, $false, "36+9=", 2) | Out-Null
$d.Content.Find.Execute("69-31=", $true, $false, $false, $false, $false, $true, 1, $false, "86-55=", 2) | Out-Null
$d.Content.Find.Execute("92-70=", $true, $false, $false, $false, $false, $true, 1, $false, "48+11=", 2) | Out-Null
$d.Content.Find.Execute("71-65=", $true, $false, $false, $false, $false, $true, 1, $false, "9+72=", 2) | Out-Null
$d.Content.Find.Execute("29+34=", $true, $false, $false, $false, $false, $true, 1, $false, "85-73=", 2) | Out-Null
$d.Content.Find.Execute("80+10=", $true, $false, $false, $false, $false, $true, 1, $false, "78+6=", 2) | Out-Null
$d.Content.Find.Execute("39-12=", $true, $false, $false, $false, $false, $true, 1, $false, "56-6=", 2) | Out-Null
$d.Content.Find.Execute("19+67=", $true, $false, $false, $false, $false, $true, 1, $false, "53+27=", 2) | Out-Null
$d.Content.Find.Execute("52+41=", $true, $false, $false, $false, $false, $true, 1, $false, "28+12=", 2) | Out-Null
$d.Content.Find.Execute("91-37=", $true, $false, $false, $false, $false, $true, 1, $false, "73-14=", 2) | Out-Null
$d.Content.Find.Execute("30+61=", $true, $false, $false, $false, $false, $true, 1, $false, "22-20=", 2) | Out-Null
$d.Content.Find.Execute("60-27=", $true, $false, $false, $false, $false, $true, 1, $false, "6+42=", 2) | Out-Null
$d.Content.Find.Execute("75-4=", $true, $false, $false, $false, $false, $true, 1, $false, "60-17=", 2) | Out-Null
$d.Content.Find.Execute("79+13=", $true, $false, $false, $false, $false, $true, 1, $false, "98-97=", 2) | Out-Null
$d.Content.Find.Execute("6+52=", $true, $false, $false, $false, $false, $true, 1, $false, "35+60=", 2) | Out-Null
$d.Content.Find.Execute("75+15=", $true, $false, $false, $false, $false, $true, 1, $false, "72+4=", 2) | Out-Null
$d.Content.Find.Execute("64+27=", $true, $false, $false, $false, $false, $true, 1, $false, "28-16=", 2) | Out-Null
$d.Content.Find.Execute("52+8=", $true, $false, $false, $false, $false, $true, 1, $false, "55+20=", 2) | Out-Null
$d.Content.Find.Execute("61+15=", $true, $false, $false, $false, $false, $true, 1, $false, "26+57=", 2) | Out-Null
$d.Content.Find.Execute("97-72=", $true, $false, $false, $false, $false, $true, 1, $false, "41-21=", 2) | Out-Null
$d.Content.Find.Execute("91-85=", $true, $false, $false, $false, $false, $true, 1, $false, "80-26=", 2) | Out-Null
$d.Content.Find.Execute("38+46=", $true, $false, $false, $false, $false, $true, 1, $false, "20+46=", 2) | Out-Null
$d.Content.Find.Execute("4+40=", $true, $false, $false, $false, $false, $true, 1, $false, "81-77=", 2) | Out-Null
$d.Content.Find.Execute("41+2=", $true, $false, $false, $false, $false, $true, 1, $false, "49-34=", 2) | Out-Null
$d.Content.Find.Execute("89+8=", $true, $false, $false, $false, $false, $true, 1, $false, "88-28=", 2) | Out-Null
$d.Content.Find.Execute("35+57=", $true, $false, $false, $false, $false, $true, 1, $false, "7+46=", 2) | Out-Null
$d.Content.Find.Execute("19+18=", $true, $false, $false, $false, $false, $true, 1, $false, "68+6=", 2) | Out-Null
$d.Content.Find.Execute("70-39=", $true, $false, $false, $false, $false, $true, 1, $false, "7+46=", 2) | Out-Null
$d.Content.Find.Execute("21+64=", $true, $false, $false, $false, $false, $true, 1, $false, "20+46=", 2) | Out-Null
$d.Content.Find.Execute("63-1=", $true, $false, $false, $false, $false, $true, 1, $false, "47-31=", 2) | Out-Null
$d.Content.Find.Execute("32-28=", $true, $false, $false, $false, $false, $true, 1, $false, "9-3=", 2) | Out-Null
$d.Content.Find.Execute("96-63=", $true, $false, $false, $false, $false, $true, 1, $false, "58-45=", 2) | Out-Null
$d.Content.Find.Execute("95-77=", $true, $false, $false, $false, $false, $true, 1, $false, "49+23=", 2) | Out-Null
$d.Content.Find.Execute("81-6=", $true, $false, $false, $false, $false, $true, 1, $false, "56-31=", 2) | Out-Null
$d.Content.Find.Execute("52-7=", $true, $false, $false, $false, $false, $true, 1, $false, "55+2=", 2) | Out-Null
$d.Content.Find.Execute("87-72=", $true, $false, $false, $false, $false, $true, 1, $false, "92-16=", 2) | Out-Null
$d.Content.Find.Execute("11+75=", $true, $false, $false, $false, $false, $true, 1, $false, "83+0=", 2) | Out-Null
$d.Content.Find.Execute("40+4=", $true, $false, $false, $false, $false, $true, 1, $false, "28+38=", 2) | Out-Null
$d.Content.Find.Execute("72-44=", $true, $false, $false, $false, $false, $true, 1, $false, "47-30=", 2) | Out-Null
$d.Content.Find.Execute("77-43=", $true, $false, $false, $false, $false, $true, 1, $false, "79+19=", 2) | Out-Null
$d.Content.Find.Execute("61-14=", $true, $false, $false, $false, $false, $true, 1, $false, "3+74=", 2) | Out-Null
$d.Content.Find.Execute("2+33=", $true, $false, $false, $false, $false, $true, 1, $false, "50-18=", 2) | Out-Null
$d.Content.Find.Execute("87-58=", $true, $false, $false, $false, $false, $true, 1, $false, "87-29=", 2) | Out-Null
$d.Content.Find.Execute("10+38=", $true, $false, $false, $false, $false, $true, 1, $false, "43+15=", 2) | Out-Null
$d.Content.Find.Execute("94-27=", $true, $false, $false, $false, $false, $true, 1, $false, "95-18=", 2) | Out-Null
$d.Content.Find.Execute("84-2=", $true, $false, $false, $false, $false, $true, 1, $false, "44-13=", 2) | Out-Null
$d.Content.Find.Execute("28+44=", $true, $false, $false, $false, $false, $true, 1, $false, "8+78=", 2) | Out-Null
$d.Content.Find.Execute("47-39=", $true, $false, $false, $false, $false, $true, 1, $false, "70+12=", 2) | Out-Null
$d.Content.Find.Execute("36+41=", $true, $false, $false, $false, $false, $true, 1, $false, "93-2=", 2) | Out-Null
$d.Content.Find.Execute("89-64=", $true, $false, $false, $false, $false, $true, 1, $false, "18+43=", 2) | Out-Null
$d.Content.Find.Execute("49-17=", $true, $false, $false, $false, $false, $true, 1, $false, "4+22=", 2) | Out-Null
$d.Content.Find.Execute("67-50=", $true, $false, $false, $false, $false, $true, 1, $false, "28+49=", 2) | Out-Null
$d.Content.Find.Execute("79-21=", $true, $false, $false, $false, $false, $true, 1, $false, "94-9=", 2) | Out-Null
$d.Content.Find.Execute("69+30=", $true, $false, $false, $false, $false, $true, 1, $false, "79-27=", 2) | Out-Null
$d.Content.Find.Execute("50+18=", $true, $false, $false, $false, $false, $true, 1, $false, "17+49=", 2) | Out-Null
$d.Content.Find.Execute("50+22=", $true, $false, $false, $false, $false, $true, 1, $false, "66+5=", 2) | Out-Null
$d.Content.Find.Execute("8+56=", $true, $false, $false, $false, $false, $true, 1, $false, "71-59=", 2) | Out-Null
$d.Content.Find.Execute("25+57=", $true, $false, $false, $false, $false, $true, 1, $false, "1+96=", 2) | Out-Null
$d.Content.Find.Execute("78-35=", $true, $false, $false, $false, $false, $true, 1, $false, "78+8=", 2) | Out-Null
$d.Content.Find.Execute("96-55=", $true, $false, $false, $false, $false, $true, 1, $false, "95-89=", 2) | Out-Null
$d.Content.Find.Execute("24+35=", $true, $false, $false, $false, $false, $true, 1, $false, "22-12=", 2) | Out-Null
$d.Content.Find.Execute("89-72=", $true, $false, $false, $false, $false, $true, 1, $false, "48-16=", 2) | Out-Null
$d.Content.Find.Execute("94-7=", $true, $false, $false, $false, $false, $true, 1, $false, "92-10=", 2) | Out-Null
$d.Content.Find.Execute("21-3=", $true, $false, $false, $false, $false, $true, 1, $false, "78-51=", 2) | Out-Null
$d.Content.Find.Execute("15+82=", $true, $false, $false, $false, $false, $true, 1, $false, "23+42=", 2) | Out-Null
$d.Content.Find.Execute("24+43=", $true, $false, $false, $false, $false, $true, 1, $false, "65-12=", 2) | Out-Null
$d.Content.Find.Execute("17-3=", $true, $false, $false, $false, $false, $true, 1, $false, "54-9=", 2) | Out-Null
$d.Content.Find.Execute("60-10=", $true, $false, $false, $false, $false, $true, 1, $false, "4+68=", 2) | Out-Null
$d.Content.Find.Execute("42-22=", $true, $false, $false, $false, $false, $true, 1, $false, "2+18=", 2) | Out-Null
$d.Content.Find.Execute("17+51=", $true, $false, $false, $false, $false, $true, 1, $false, "96-69=", 2) | Out-Null
$d.Content.Find.Execute("23+68=", $true, $false, $false, $false, $false, $true, 1, $false, "26+1=", 2) | Out-Null
$d.Content.Find.Execute("89-38=", $true, $false, $false, $false, $false, $true, 1, $false, "13+51=", 2) | Out-Null
$d.Content.Find.Execute("0+39=", $true, $false, $false, $false, $false, $true, 1, $false, "59+37=", 2) | Out-Null
$d.Content.Find.Execute("27+60=", $true, $false, $false, $false, $false, $true, 1, $false, "88-12=", 2) | Out-Null
$d.Content.Find.Execute("70-62=", $true, $false, $false, $false, $false, $true, 1, $false, "94-68=", 2) | Out-Null
$d.Content.Find.Execute("91-50=", $true, $false, $false, $false, $false, $true, 1, $false, "12-12=", 2) | Out-Null
$d.Content.Find.Execute("3+20=", $true, $false, $false, $false, $false, $true, 1, $false, "74-30=", 2) | Out-Null
$d.Content.Find.Execute("29+22=", $true, $false, $false, $false, $false, $true, 1, $false, "3+23=", 2) | Out-Null
